# Updated symbol list on Wed Dec 21 07:35:51 UTC 2022 with GitHub Actions
#
# This script updates the "Price" column (D) and a couple of "Volume(1h)"
# column (E) text labels on the active worksheet to reflect refreshed
# crypto market data. All of the target cells hold numeric-looking text
# (t="inlineStr"/shared-string, not real numbers), so each cell's number
# format is forced to Text ("@") before the value is assigned, and then
# restored to General afterwards -- this keeps Excel from silently
# re-interpreting strings like "248.24" as a floating point number (which
# would introduce binary rounding noise such as 248.24000000000001) while
# leaving the cell's effective formatting unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Value
    )
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Value
    $cell.NumberFormat = "General"
}

Set-TextValue "D2"  "248.24"
Set-TextValue "D3"  "22.63"
Set-TextValue "D4"  "5.389"
Set-TextValue "D5"  "0.05688"
Set-TextValue "D7"  "6.320"
Set-TextValue "D8"  "0.8129"
Set-TextValue "D9"  "0.9273"
Set-TextValue "D10" "0.1412"
Set-TextValue "D11" "0.07453"
Set-TextValue "D12" "0.03079"
Set-TextValue "D13" "0.03016"
Set-TextValue "D14" "0.09376"
Set-TextValue "D15" "3.772"
Set-TextValue "D16" "0.001585"
Set-TextValue "D17" "0.04758"

$ws.Range("E19").Value = "18OneONE"

Set-TextValue "D20" "0.006459"
Set-TextValue "D21" "0.004995"
Set-TextValue "D22" "0.001026"
Set-TextValue "D24" "3.699"
Set-TextValue "D25" "2.153"
Set-TextValue "D40" "0.03995"
Set-TextValue "D41" "0.006830"
Set-TextValue "D42" "0.1067"
Set-TextValue "D43" "0.002710"
Set-TextValue "D44" "0.007476"
Set-TextValue "D45" "0.00005801"

Set-TextValue "D47" "0.4300"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"

Set-TextValue "D48" "0.2150"
Set-TextValue "D49" "0.00002100"
